$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot (column D holds numeric-looking
# values that are stored as TEXT in this workbook, e.g. "0.006410" or
# "0.00000000750", where trailing zeros and decimal formatting are significant).
# Each price cell is set to Text format before its value is written so the
# new value is kept as an exact string instead of being parsed into a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "264.11"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.65"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.274"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06166"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.577"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.543"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.394"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8239"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1624"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08215"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03545"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03187"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09212"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.778"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001623"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04669"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006412"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006171"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001070"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001501"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.724"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.235"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01356"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1244"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002719"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04690"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006996"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003761"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1117"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01194"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006120"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009920"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9822"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00001902"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01241"
